# Advertising_Campaign.xlsx - translate Dutch labels to English
# (per commit: "Topic colofon translated. ... Corrections to practice
#  file Advertising campaign.")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Number of flyers sent"
$ws.Range("A2").Value = "Response (%)"
$ws.Range("A4").Value = "Printing costs per piece"
$ws.Range("A5").Value = "Shipping costs per piece"
$ws.Range("A6").Value = "Revenue per response"
$ws.Range("A8").Value = "Response (number)"
$ws.Range("A9").Value = "Revenue total"
$ws.Range("A10").Value = "Costs total"
$ws.Range("A11").Value = "profit"

$ws.Range("A1").Select()
